$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About": rework the Sources / Notes section with the new EPA citation
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("About")

# The old "last updated" date stamp in C1 is gone entirely (column C is no
# longer used on this sheet).
$ws.Range("C1").Clear()

# Make room for the new 4-line source citation (year / author / title / url)
# right under the existing "Sources:" row.
$ws.Range("A4:A7").EntireRow.Insert()

# Make room for the new 2-line explanatory note about the chosen values,
# inserted right before the "For more on this..." explanation.
$ws.Range("A16:A18").EntireRow.Insert()

# Replace the old "None needed.  Handled through calibration." source note
# with the new EPA citation (filled in the same order the original author
# typed them, so the shared-string table comes out in the same order).
$ws.Range("B7").Value = "Table 5 Generalized Cost Coefficient Calibration"
$ws.Range("B3").Value = "United States EPA"
$ws.Range("B5").Value = "Consumer Vehicle Choice Model Documentation"
$ws.Range("B6").Value = "https://nepis.epa.gov/Exe/ZyPDF.cgi/P100EZ37.PDF?Dockey=P100EZ37.PDF"

$ws.Range("B4").Value = 2012
$ws.Range("B4").HorizontalAlignment = -4131

$ws.Range("A16").Value = "We choose a value of -3 for passenger vehicles and a value of -5 for other vehicle types, "
$ws.Range("A17").Value = "based on the ranges in Table 5 of the cited EPA documentation."

# ---------------------------------------------------------------------------
# Sheet "TTLE": update the logit exponent values from -3 to -5
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("TTLE")
$ws2.Range("B2:C7").Value = -5

Write-Output "edit applied"
